$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1), copying the formatting used by the other
# header cells (bold font, thin border, centered/top alignment - same as G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell for the Save column (H2)
$ws.Range("H2").Value = 0
